$wb = $excel.ActiveWorkbook

# --- Sheet "Estadisticos 2P": segundo parcial grades entered, recompute
#     Blancos (D), Reprobados (E), Aprobados (F), Por_Apro (G) and add Promedio (H) ---
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 0
$ws2P.Range("F2").Value = 34
$ws2P.Range("G2").Value = 100
$ws2P.Range("H2").Value = 8.3

$ws2P.Range("D3").Value = 0
$ws2P.Range("E3").Value = 1
$ws2P.Range("F3").Value = 40
$ws2P.Range("G3").Value = 97.56
$ws2P.Range("H3").Value = 8.7

$ws2P.Range("D4").Value = 0
$ws2P.Range("E4").Value = 1
$ws2P.Range("F4").Value = 40
$ws2P.Range("G4").Value = 97.56
$ws2P.Range("H4").Value = 8.7

$ws2P.Range("D5").Value = 0
$ws2P.Range("E5").Value = 0
$ws2P.Range("F5").Value = 36
$ws2P.Range("G5").Value = 100
$ws2P.Range("H5").Value = 9

$ws2P.Range("D6").Value = 0
$ws2P.Range("E6").Value = 0
$ws2P.Range("F6").Value = 36
$ws2P.Range("G6").Value = 100
$ws2P.Range("H6").Value = 9

# --- Sheet "Estadisticos Final": recompute Promedio (H) now that 2P is in ---
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("H2").Value = 9.2
$wsFinal.Range("H3").Value = 9
$wsFinal.Range("H4").Value = 9
$wsFinal.Range("H5").Value = 9.3
$wsFinal.Range("H6").Value = 9.2

# --- Sheet "Rescatables": update remaining exams (G) count ---
$wsResc = $wb.Worksheets.Item("Rescatables")

$wsResc.Range("G2").Value = 3
$wsResc.Range("G3").Value = 3
